$d = $word.ActiveDocument

$rPr32 = '<w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'
$p1Xml = '<w:p w14:paraId="370D496E" w14:textId="270058FF" w:rsidR="00F31009" w:rsidRDefault="00B02BE4"><w:pPr>' + $rPr32 + '</w:pPr><w:r>' + $rPr32 + '<w:t xml:space="preserve">Sprint </w:t></w:r><w:proofErr w:type="spellStart"/><w:r>' + $rPr32 + '<w:t>zero</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r>' + $rPr32 + '<w:t xml:space="preserve"> - aloituspalaveri 9.11.2023</w:t></w:r></w:p>'

$pkgTemplate = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">{0}</w:document></pkg:xmlData></pkg:part></pkg:package>'

$body1 = '<w:body>' + $p1Xml + '</w:body>'
$xml1 = $pkgTemplate -f $body1

$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML($xml1)

$rPr28 = '<w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$p2Xml = '<w:p w14:paraId="44B865D5" w14:textId="780B8997" w:rsidR="00DC2B2F" w:rsidRDefault="00DC2B2F" w:rsidP="00DC2B2F"><w:pPr>' + $rPr28 + '</w:pPr><w:r>' + $rPr28 + '<w:t>Sprint 0 Retrospektiivi:</w:t></w:r></w:p>'
$body2 = '<w:body>' + $p2Xml + '</w:body>'
$xml2 = $pkgTemplate -f $body2

$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Sprint 0 Retrospektiivi:`r") {
        $p2 = $p
        break
    }
}
$p2.Range.InsertXML($xml2)
